$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.680.90"
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").Value = "'1.920.51"
$ws.Range("E3").Value = '  +1.55%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = "'240.11"
$ws.Range("E5").Value = '  -2.09%  '

$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").Value = "'0.4943"
$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").Value = "'0.3002"
$ws.Range("E8").Value = '  +1.32%  '

$ws.Range("D9").Value = "'0.06771"
$ws.Range("E9").Value = '  -0.38%  '

$ws.Range("D10").Value = "'1.936.92"
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("D11").Value = "'17.22"
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("D12").Value = "'0.07356"
$ws.Range("E12").Value = '  +1.75%  '

$ws.Range("D13").Value = "'5.213"
$ws.Range("E13").Value = '  +3.10%  '

$ws.Range("D14").Value = "'88.72"
$ws.Range("E14").Value = '  -2.89%  '

$ws.Range("D15").Value = "'0.6748"
$ws.Range("E15").Value = '  -0.58%  '

$ws.Range("D16").Value = "'30.656.35"
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").Value = "'0.000007964"
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").Value = "'13.56"
$ws.Range("E18").Value = '  +2.66%  '

$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").Value = "'2.150.64"
$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("D21").Value = "'5.419"
$ws.Range("E21").Value = '  +12.29%  '

$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("D23").Value = "'197.68"
$ws.Range("E23").Value = '  +1.78%  '

$ws.Range("D24").Value = "'6.346"
$ws.Range("E24").Value = '  +3.77%  '

$ws.Range("D25").Value = "'9.667"
$ws.Range("E25").Value = '  +3.15%  '

$ws.Range("D26").Value = "'164.27"
$ws.Range("E26").Value = '  +5.89%  '

$ws.Range("D27").Value = "'18.70"
$ws.Range("E27").Value = '  -3.17%  '

$ws.Range("D28").Value = "'1.962"
$ws.Range("E28").Value = '  +2.83%  '

$ws.Range("D29").Value = "'1.474"
$ws.Range("E29").Value = '  +4.84%  '

$ws.Range("D30").Value = "'4.377"
$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("D31").Value = "'0.09172"
$ws.Range("E31").Value = '  +0.74%  '

$ws.Range("D32").Value = "'4.083"
$ws.Range("E32").Value = '  +1.61%  '

$ws.Range("D33").Value = "'0.05273"
$ws.Range("E33").Value = '  +1.30%  '

$ws.Range("D34").Value = "'0.7433"
$ws.Range("E34").Value = '  -2.77%  '

$ws.Range("D35").Value = "'1.119"
$ws.Range("E35").Value = '  +0.57%  '

$ws.Range("D36").Value = "'2.717"
$ws.Range("E36").Value = '  -2.04%  '

$ws.Range("D37").Value = "'0.01846"
$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("D38").Value = "'2.722"
$ws.Range("E38").Value = '  +1.58%  '

$ws.Range("D39").Value = "'0.9283"
$ws.Range("E39").Value = '  -0.86%  '

$ws.Range("D40").Value = "'2.091"
$ws.Range("E40").Value = '  -2.84%  '

$ws.Range("D41").Value = "'0.4484"
$ws.Range("E41").Value = '  +1.02%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'5.967"
$ws.Range("E42").Value = '  +3.36%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = "'71.91"
$ws.Range("E43").Value = '  +24.34%  '

$ws.Range("D44").Value = "'106.46"
$ws.Range("E44").Value = '  +0.73%  '

$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = '  +0.27%  '

$ws.Range("D46").Value = "'0.1401"
$ws.Range("E46").Value = '  +4.01%  '

$ws.Range("D47").Value = "'7.676"
$ws.Range("E47").Value = '  +0.64%  '

$ws.Range("D48").Value = "'9.101"
$ws.Range("E48").Value = '  +4.41%  '

$ws.Range("D49").Value = "'35.24"
$ws.Range("E49").Value = '  +4.51%  '

$ws.Range("D50").Value = "'0.05893"
$ws.Range("E50").Value = '  +0.48%  '

$ws.Range("D51").Value = "'0.4048"
$ws.Range("E51").Value = '  +2.70%  '
